# Refresh the "Saldo" export: drop stale accounts, add new ones, update balances,
# and re-sort the data rows by balance (descending) to match the refreshed report.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The refreshed export has 6 fewer data rows than before; delete that many rows from
# the data block so the trailing blank separator + "Filtros aplicados" footer rows land
# back in the correct position once every row below the header is rewritten.
$ws.Range("A216:A221").EntireRow.Delete()

# Keep the "Conta" column as text so account numbers keep their leading zeros.
$ws.Range("A2:A215").NumberFormat = "@"

$ws.Cells.Item(2, 1).Value = "005305448"
$ws.Cells.Item(2, 2).Value = "ALPHASITIO"
$ws.Cells.Item(2, 3).Value = 365703.31
$ws.Cells.Item(3, 1).Value = "004479287"
$ws.Cells.Item(3, 2).Value = "ANA"
$ws.Cells.Item(3, 3).Value = 65727.36
$ws.Cells.Item(4, 1).Value = "005070742"
$ws.Cells.Item(4, 2).Value = "JUSCELINO"
$ws.Cells.Item(4, 3).Value = 64399.9
$ws.Cells.Item(5, 1).Value = "004935287"
$ws.Cells.Item(5, 2).Value = "ODILON"
$ws.Cells.Item(5, 3).Value = 30297.77
$ws.Cells.Item(6, 1).Value = "004224011"
$ws.Cells.Item(6, 2).Value = "THOMAS"
$ws.Cells.Item(6, 3).Value = 21582.35
$ws.Cells.Item(7, 1).Value = "004454365"
$ws.Cells.Item(7, 2).Value = "RAFAEL"
$ws.Cells.Item(7, 3).Value = 13735.23
$ws.Cells.Item(8, 1).Value = "004361159"
$ws.Cells.Item(8, 2).Value = "HFR"
$ws.Cells.Item(8, 3).Value = 5714.31
$ws.Cells.Item(9, 1).Value = "004229526"
$ws.Cells.Item(9, 2).Value = "EDUARDA"
$ws.Cells.Item(9, 3).Value = 5101.25
$ws.Cells.Item(10, 1).Value = "004643153"
$ws.Cells.Item(10, 2).Value = "CARLA"
$ws.Cells.Item(10, 3).Value = 1024.32
$ws.Cells.Item(11, 1).Value = "004488571"
$ws.Cells.Item(11, 2).Value = "CARLOS"
$ws.Cells.Item(11, 3).Value = 1000
$ws.Cells.Item(12, 1).Value = "004504449"
$ws.Cells.Item(12, 2).Value = "KELMA"
$ws.Cells.Item(12, 3).Value = 1000
$ws.Cells.Item(13, 1).Value = "004240014"
$ws.Cells.Item(13, 2).Value = "ISABELE"
$ws.Cells.Item(13, 3).Value = 998.71
$ws.Cells.Item(14, 1).Value = "004334158"
$ws.Cells.Item(14, 2).Value = "LEONE"
$ws.Cells.Item(14, 3).Value = 994.66
$ws.Cells.Item(15, 1).Value = "001882235"
$ws.Cells.Item(15, 2).Value = "LAGO"
$ws.Cells.Item(15, 3).Value = 966.83
$ws.Cells.Item(16, 1).Value = "004457389"
$ws.Cells.Item(16, 2).Value = "RAFAEL"
$ws.Cells.Item(16, 3).Value = 926.76
$ws.Cells.Item(17, 1).Value = "004487016"
$ws.Cells.Item(17, 2).Value = "ROGERIO"
$ws.Cells.Item(17, 3).Value = 921.71
$ws.Cells.Item(18, 1).Value = "004392159"
$ws.Cells.Item(18, 2).Value = "RODRIGO"
$ws.Cells.Item(18, 3).Value = 900.21
$ws.Cells.Item(19, 1).Value = "004855570"
$ws.Cells.Item(19, 2).Value = "LUISA"
$ws.Cells.Item(19, 3).Value = 895.15
$ws.Cells.Item(20, 1).Value = "004550605"
$ws.Cells.Item(20, 2).Value = "REJANE"
$ws.Cells.Item(20, 3).Value = 869.47
$ws.Cells.Item(21, 1).Value = "004936634"
$ws.Cells.Item(21, 2).Value = "LEONARDO"
$ws.Cells.Item(21, 3).Value = 865.79
$ws.Cells.Item(22, 1).Value = "004211368"
$ws.Cells.Item(22, 2).Value = "ILTON"
$ws.Cells.Item(22, 3).Value = 856.39
$ws.Cells.Item(23, 1).Value = "004216504"
$ws.Cells.Item(23, 2).Value = "WANDER"
$ws.Cells.Item(23, 3).Value = 850.24
$ws.Cells.Item(24, 1).Value = "004563252"
$ws.Cells.Item(24, 2).Value = "FERNANDO"
$ws.Cells.Item(24, 3).Value = 837.81
$ws.Cells.Item(25, 1).Value = "005245032"
$ws.Cells.Item(25, 2).Value = "ROSA"
$ws.Cells.Item(25, 3).Value = 824.46
$ws.Cells.Item(26, 1).Value = "002064834"
$ws.Cells.Item(26, 2).Value = "RAFAELA"
$ws.Cells.Item(26, 3).Value = 813.54
$ws.Cells.Item(27, 1).Value = "004452476"
$ws.Cells.Item(27, 2).Value = "IVONE"
$ws.Cells.Item(27, 3).Value = 768.76
$ws.Cells.Item(28, 1).Value = "004381180"
$ws.Cells.Item(28, 2).Value = "HFR"
$ws.Cells.Item(28, 3).Value = 743.31
$ws.Cells.Item(29, 1).Value = "002697806"
$ws.Cells.Item(29, 2).Value = "CLAUDIA"
$ws.Cells.Item(29, 3).Value = 705.53
$ws.Cells.Item(30, 1).Value = "004388077"
$ws.Cells.Item(30, 2).Value = "WLADMIR"
$ws.Cells.Item(30, 3).Value = 703.17
$ws.Cells.Item(31, 1).Value = "004359408"
$ws.Cells.Item(31, 2).Value = "HEPTA"
$ws.Cells.Item(31, 3).Value = 693.95
$ws.Cells.Item(32, 1).Value = "004646727"
$ws.Cells.Item(32, 2).Value = "RENATA"
$ws.Cells.Item(32, 3).Value = 679.08
$ws.Cells.Item(33, 1).Value = "004975924"
$ws.Cells.Item(33, 2).Value = "SERGIO"
$ws.Cells.Item(33, 3).Value = 672.05
$ws.Cells.Item(34, 1).Value = "004972070"
$ws.Cells.Item(34, 2).Value = "MARIA"
$ws.Cells.Item(34, 3).Value = 670.51
$ws.Cells.Item(35, 1).Value = "004252768"
$ws.Cells.Item(35, 2).Value = "ALESSANDRO"
$ws.Cells.Item(35, 3).Value = 656.22
$ws.Cells.Item(36, 1).Value = "004481463"
$ws.Cells.Item(36, 2).Value = "MARA"
$ws.Cells.Item(36, 3).Value = 637.03
$ws.Cells.Item(37, 1).Value = "004228456"
$ws.Cells.Item(37, 2).Value = "FLASH"
$ws.Cells.Item(37, 3).Value = 611.4
$ws.Cells.Item(38, 1).Value = "004517080"
$ws.Cells.Item(38, 2).Value = "TATIANA"
$ws.Cells.Item(38, 3).Value = 607.94
$ws.Cells.Item(39, 1).Value = "005079311"
$ws.Cells.Item(39, 2).Value = "JOVINO"
$ws.Cells.Item(39, 3).Value = 603.35
$ws.Cells.Item(40, 1).Value = "004574428"
$ws.Cells.Item(40, 2).Value = "GUILHERME"
$ws.Cells.Item(40, 3).Value = 596.31
$ws.Cells.Item(41, 1).Value = "005046919"
$ws.Cells.Item(41, 2).Value = "MARIANA"
$ws.Cells.Item(41, 3).Value = 590.31
$ws.Cells.Item(42, 1).Value = "004806244"
$ws.Cells.Item(42, 2).Value = "CARLA"
$ws.Cells.Item(42, 3).Value = 585.78
$ws.Cells.Item(43, 1).Value = "005142592"
$ws.Cells.Item(43, 2).Value = "ALBERTO"
$ws.Cells.Item(43, 3).Value = 551.34
$ws.Cells.Item(44, 1).Value = "004238436"
$ws.Cells.Item(44, 2).Value = "DIEGO"
$ws.Cells.Item(44, 3).Value = 547.62
$ws.Cells.Item(45, 1).Value = "005055865"
$ws.Cells.Item(45, 2).Value = "G3C"
$ws.Cells.Item(45, 3).Value = 526.02
$ws.Cells.Item(46, 1).Value = "004395314"
$ws.Cells.Item(46, 2).Value = "MARIA"
$ws.Cells.Item(46, 3).Value = 522.36
$ws.Cells.Item(47, 1).Value = "004398174"
$ws.Cells.Item(47, 2).Value = "DANIELE"
$ws.Cells.Item(47, 3).Value = 496.52
$ws.Cells.Item(48, 1).Value = "004322719"
$ws.Cells.Item(48, 2).Value = "GISELA"
$ws.Cells.Item(48, 3).Value = 493.87
$ws.Cells.Item(49, 1).Value = "004448303"
$ws.Cells.Item(49, 2).Value = "NASSIM"
$ws.Cells.Item(49, 3).Value = 488.3
$ws.Cells.Item(50, 1).Value = "004556853"
$ws.Cells.Item(50, 2).Value = "MARCEL"
$ws.Cells.Item(50, 3).Value = 483.1
$ws.Cells.Item(51, 1).Value = "004387250"
$ws.Cells.Item(51, 2).Value = "MONICA"
$ws.Cells.Item(51, 3).Value = 478
$ws.Cells.Item(52, 1).Value = "005338054"
$ws.Cells.Item(52, 2).Value = "ELAINE"
$ws.Cells.Item(52, 3).Value = 475.92
$ws.Cells.Item(53, 1).Value = "000772433"
$ws.Cells.Item(53, 2).Value = "MARCELO"
$ws.Cells.Item(53, 3).Value = 465.49
$ws.Cells.Item(54, 1).Value = "005266369"
$ws.Cells.Item(54, 2).Value = "EG"
$ws.Cells.Item(54, 3).Value = 459.39
$ws.Cells.Item(55, 1).Value = "005203562"
$ws.Cells.Item(55, 2).Value = "ROBERIO"
$ws.Cells.Item(55, 3).Value = 453.65
$ws.Cells.Item(56, 1).Value = "004893911"
$ws.Cells.Item(56, 2).Value = "FAUSTO"
$ws.Cells.Item(56, 3).Value = 449.88
$ws.Cells.Item(57, 1).Value = "004415557"
$ws.Cells.Item(57, 2).Value = "FILIPE"
$ws.Cells.Item(57, 3).Value = 427.15
$ws.Cells.Item(58, 1).Value = "004276856"
$ws.Cells.Item(58, 2).Value = "DAURA"
$ws.Cells.Item(58, 3).Value = 426.61
$ws.Cells.Item(59, 1).Value = "004556150"
$ws.Cells.Item(59, 2).Value = "MARINA"
$ws.Cells.Item(59, 3).Value = 409.35
$ws.Cells.Item(60, 1).Value = "004480970"
$ws.Cells.Item(60, 2).Value = "ALBERTO"
$ws.Cells.Item(60, 3).Value = 399.15
$ws.Cells.Item(61, 1).Value = "004756981"
$ws.Cells.Item(61, 2).Value = "MATEUS"
$ws.Cells.Item(61, 3).Value = 357.49
$ws.Cells.Item(62, 1).Value = "004587511"
$ws.Cells.Item(62, 2).Value = "CARLOS"
$ws.Cells.Item(62, 3).Value = 352.61
$ws.Cells.Item(63, 1).Value = "004289402"
$ws.Cells.Item(63, 2).Value = "LARISSA"
$ws.Cells.Item(63, 3).Value = 349.86
$ws.Cells.Item(64, 1).Value = "004381415"
$ws.Cells.Item(64, 2).Value = "JOAO"
$ws.Cells.Item(64, 3).Value = 349.74
$ws.Cells.Item(65, 1).Value = "004212476"
$ws.Cells.Item(65, 2).Value = "MARIA"
$ws.Cells.Item(65, 3).Value = 349.25
$ws.Cells.Item(66, 1).Value = "004204500"
$ws.Cells.Item(66, 2).Value = "EDWARD"
$ws.Cells.Item(66, 3).Value = 343.71
$ws.Cells.Item(67, 1).Value = "005009992"
$ws.Cells.Item(67, 2).Value = "ALINE"
$ws.Cells.Item(67, 3).Value = 330.17
$ws.Cells.Item(68, 1).Value = "005101676"
$ws.Cells.Item(68, 2).Value = "ELENI"
$ws.Cells.Item(68, 3).Value = 314.83
$ws.Cells.Item(69, 1).Value = "004214592"
$ws.Cells.Item(69, 2).Value = "MERG"
$ws.Cells.Item(69, 3).Value = 312.75
$ws.Cells.Item(70, 1).Value = "004480134"
$ws.Cells.Item(70, 2).Value = "JOSE"
$ws.Cells.Item(70, 3).Value = 308.81
$ws.Cells.Item(71, 1).Value = "004332103"
$ws.Cells.Item(71, 2).Value = "JOSE"
$ws.Cells.Item(71, 3).Value = 300.31
$ws.Cells.Item(72, 1).Value = "005121919"
$ws.Cells.Item(72, 2).Value = "JORGE"
$ws.Cells.Item(72, 3).Value = 297.95
$ws.Cells.Item(73, 1).Value = "004259659"
$ws.Cells.Item(73, 2).Value = "BENTO"
$ws.Cells.Item(73, 3).Value = 293.95
$ws.Cells.Item(74, 1).Value = "000330949"
$ws.Cells.Item(74, 2).Value = "RENATO"
$ws.Cells.Item(74, 3).Value = 285.36
$ws.Cells.Item(75, 1).Value = "000626491"
$ws.Cells.Item(75, 2).Value = "FELIPE"
$ws.Cells.Item(75, 3).Value = 280.06
$ws.Cells.Item(76, 1).Value = "004355790"
$ws.Cells.Item(76, 2).Value = "MINEIA"
$ws.Cells.Item(76, 3).Value = 279.74
$ws.Cells.Item(77, 1).Value = "004486497"
$ws.Cells.Item(77, 2).Value = "ELENA"
$ws.Cells.Item(77, 3).Value = 257.45
$ws.Cells.Item(78, 1).Value = "004927044"
$ws.Cells.Item(78, 2).Value = "CINTIA"
$ws.Cells.Item(78, 3).Value = 257.37
$ws.Cells.Item(79, 1).Value = "004580355"
$ws.Cells.Item(79, 2).Value = "LARISSA"
$ws.Cells.Item(79, 3).Value = 227.65
$ws.Cells.Item(80, 1).Value = "004870976"
$ws.Cells.Item(80, 2).Value = "HFR"
$ws.Cells.Item(80, 3).Value = 222.74
$ws.Cells.Item(81, 1).Value = "004526450"
$ws.Cells.Item(81, 2).Value = "MSD"
$ws.Cells.Item(81, 3).Value = 205.85
$ws.Cells.Item(82, 1).Value = "004334062"
$ws.Cells.Item(82, 2).Value = "MERG"
$ws.Cells.Item(82, 3).Value = 205.56
$ws.Cells.Item(83, 1).Value = "004475395"
$ws.Cells.Item(83, 2).Value = "DAVID"
$ws.Cells.Item(83, 3).Value = 185.02
$ws.Cells.Item(84, 1).Value = "004328934"
$ws.Cells.Item(84, 2).Value = "VALERIA"
$ws.Cells.Item(84, 3).Value = 182.14
$ws.Cells.Item(85, 1).Value = "004360431"
$ws.Cells.Item(85, 2).Value = "CARLOS"
$ws.Cells.Item(85, 3).Value = 164.01
$ws.Cells.Item(86, 1).Value = "004511696"
$ws.Cells.Item(86, 2).Value = "KRYSCIA"
$ws.Cells.Item(86, 3).Value = 150.47
$ws.Cells.Item(87, 1).Value = "005022526"
$ws.Cells.Item(87, 2).Value = "ALEXANDRE"
$ws.Cells.Item(87, 3).Value = 147.18
$ws.Cells.Item(88, 1).Value = "004493324"
$ws.Cells.Item(88, 2).Value = "DANIEL"
$ws.Cells.Item(88, 3).Value = 143.45
$ws.Cells.Item(89, 1).Value = "005141215"
$ws.Cells.Item(89, 2).Value = "KARINA"
$ws.Cells.Item(89, 3).Value = 137.66
$ws.Cells.Item(90, 1).Value = "004243043"
$ws.Cells.Item(90, 2).Value = "SUELI"
$ws.Cells.Item(90, 3).Value = 134.8
$ws.Cells.Item(91, 1).Value = "005274028"
$ws.Cells.Item(91, 2).Value = "RAFAEL"
$ws.Cells.Item(91, 3).Value = 132.43
$ws.Cells.Item(92, 1).Value = "004435987"
$ws.Cells.Item(92, 2).Value = "MARCO"
$ws.Cells.Item(92, 3).Value = 125.33
$ws.Cells.Item(93, 1).Value = "004211911"
$ws.Cells.Item(93, 2).Value = "ZENILDA"
$ws.Cells.Item(93, 3).Value = 120
$ws.Cells.Item(94, 1).Value = "004404342"
$ws.Cells.Item(94, 2).Value = "ADSON"
$ws.Cells.Item(94, 3).Value = 115.85
$ws.Cells.Item(95, 1).Value = "004754920"
$ws.Cells.Item(95, 2).Value = "LUIS"
$ws.Cells.Item(95, 3).Value = 114.69
$ws.Cells.Item(96, 1).Value = "004421636"
$ws.Cells.Item(96, 2).Value = "PATRICIA"
$ws.Cells.Item(96, 3).Value = 110
$ws.Cells.Item(97, 1).Value = "004221638"
$ws.Cells.Item(97, 2).Value = "CAROLINE"
$ws.Cells.Item(97, 3).Value = 109.24
$ws.Cells.Item(98, 1).Value = "004536602"
$ws.Cells.Item(98, 2).Value = "TATIANY"
$ws.Cells.Item(98, 3).Value = 108.62
$ws.Cells.Item(99, 1).Value = "002687737"
$ws.Cells.Item(99, 2).Value = "JOSE"
$ws.Cells.Item(99, 3).Value = 101.02
$ws.Cells.Item(100, 1).Value = "005040864"
$ws.Cells.Item(100, 2).Value = "ANDRE"
$ws.Cells.Item(100, 3).Value = 100
$ws.Cells.Item(101, 1).Value = "004908680"
$ws.Cells.Item(101, 2).Value = "ELENE"
$ws.Cells.Item(101, 3).Value = 99.31
$ws.Cells.Item(102, 1).Value = "004472076"
$ws.Cells.Item(102, 2).Value = "RUBENS"
$ws.Cells.Item(102, 3).Value = 99.18
$ws.Cells.Item(103, 1).Value = "004339183"
$ws.Cells.Item(103, 2).Value = "JALISON"
$ws.Cells.Item(103, 3).Value = 95.69
$ws.Cells.Item(104, 1).Value = "004431591"
$ws.Cells.Item(104, 2).Value = "MARIO"
$ws.Cells.Item(104, 3).Value = 93.87
$ws.Cells.Item(105, 1).Value = "005256849"
$ws.Cells.Item(105, 2).Value = "SANDRO"
$ws.Cells.Item(105, 3).Value = 92.78
$ws.Cells.Item(106, 1).Value = "004335031"
$ws.Cells.Item(106, 2).Value = "EDMUNDO"
$ws.Cells.Item(106, 3).Value = 92.73
$ws.Cells.Item(107, 1).Value = "004350197"
$ws.Cells.Item(107, 2).Value = "GISELA"
$ws.Cells.Item(107, 3).Value = 91.94
$ws.Cells.Item(108, 1).Value = "004239387"
$ws.Cells.Item(108, 2).Value = "LUIZ"
$ws.Cells.Item(108, 3).Value = 89.82
$ws.Cells.Item(109, 1).Value = "004212132"
$ws.Cells.Item(109, 2).Value = "JOAO"
$ws.Cells.Item(109, 3).Value = 86.38
$ws.Cells.Item(110, 1).Value = "004207374"
$ws.Cells.Item(110, 2).Value = "ANGELICA"
$ws.Cells.Item(110, 3).Value = 85.13
$ws.Cells.Item(111, 1).Value = "005035754"
$ws.Cells.Item(111, 2).Value = "JOSE"
$ws.Cells.Item(111, 3).Value = 83.31
$ws.Cells.Item(112, 1).Value = "004216657"
$ws.Cells.Item(112, 2).Value = "JOAO"
$ws.Cells.Item(112, 3).Value = 80.63
$ws.Cells.Item(113, 1).Value = "004318604"
$ws.Cells.Item(113, 2).Value = "RENAN"
$ws.Cells.Item(113, 3).Value = 80.51
$ws.Cells.Item(114, 1).Value = "004451996"
$ws.Cells.Item(114, 2).Value = "ADRIANO"
$ws.Cells.Item(114, 3).Value = 80.36
$ws.Cells.Item(115, 1).Value = "004267976"
$ws.Cells.Item(115, 2).Value = "E3"
$ws.Cells.Item(115, 3).Value = 79.84
$ws.Cells.Item(116, 1).Value = "001294033"
$ws.Cells.Item(116, 2).Value = "VIVIANE"
$ws.Cells.Item(116, 3).Value = 79.82
$ws.Cells.Item(117, 1).Value = "005009922"
$ws.Cells.Item(117, 2).Value = "ANA"
$ws.Cells.Item(117, 3).Value = 79.02
$ws.Cells.Item(118, 1).Value = "004470679"
$ws.Cells.Item(118, 2).Value = "RODOLFO"
$ws.Cells.Item(118, 3).Value = 77.51
$ws.Cells.Item(119, 1).Value = "003115072"
$ws.Cells.Item(119, 2).Value = "VICTOR"
$ws.Cells.Item(119, 3).Value = 69.1
$ws.Cells.Item(120, 1).Value = "005133039"
$ws.Cells.Item(120, 2).Value = "PAULO"
$ws.Cells.Item(120, 3).Value = 66.51
$ws.Cells.Item(121, 1).Value = "004855596"
$ws.Cells.Item(121, 2).Value = "MARIANA"
$ws.Cells.Item(121, 3).Value = 64.36
$ws.Cells.Item(122, 1).Value = "004335251"
$ws.Cells.Item(122, 2).Value = "EDMUNDO"
$ws.Cells.Item(122, 3).Value = 62.39
$ws.Cells.Item(123, 1).Value = "000834301"
$ws.Cells.Item(123, 2).Value = "MARCUS"
$ws.Cells.Item(123, 3).Value = 57.13
$ws.Cells.Item(124, 1).Value = "004588677"
$ws.Cells.Item(124, 2).Value = "RACHEL"
$ws.Cells.Item(124, 3).Value = 55.91
$ws.Cells.Item(125, 1).Value = "004517506"
$ws.Cells.Item(125, 2).Value = "LUIZ"
$ws.Cells.Item(125, 3).Value = 55.87
$ws.Cells.Item(126, 1).Value = "004215217"
$ws.Cells.Item(126, 2).Value = "CAROLINA"
$ws.Cells.Item(126, 3).Value = 55.66
$ws.Cells.Item(127, 1).Value = "004321092"
$ws.Cells.Item(127, 2).Value = "DANIEL"
$ws.Cells.Item(127, 3).Value = 55.23
$ws.Cells.Item(128, 1).Value = "004329229"
$ws.Cells.Item(128, 2).Value = "GABRIEL"
$ws.Cells.Item(128, 3).Value = 54.91
$ws.Cells.Item(129, 1).Value = "005032151"
$ws.Cells.Item(129, 2).Value = "ANA"
$ws.Cells.Item(129, 3).Value = 52.9
$ws.Cells.Item(130, 1).Value = "004268684"
$ws.Cells.Item(130, 2).Value = "PATRICIA"
$ws.Cells.Item(130, 3).Value = 52.7
$ws.Cells.Item(131, 1).Value = "004400640"
$ws.Cells.Item(131, 2).Value = "FELIPE"
$ws.Cells.Item(131, 3).Value = 51.44
$ws.Cells.Item(132, 1).Value = "004115403"
$ws.Cells.Item(132, 2).Value = "HEBERT"
$ws.Cells.Item(132, 3).Value = 50.87
$ws.Cells.Item(133, 1).Value = "004208447"
$ws.Cells.Item(133, 2).Value = "LEILA"
$ws.Cells.Item(133, 3).Value = 50
$ws.Cells.Item(134, 1).Value = "004278033"
$ws.Cells.Item(134, 2).Value = "DAISY"
$ws.Cells.Item(134, 3).Value = 47.37
$ws.Cells.Item(135, 1).Value = "005216881"
$ws.Cells.Item(135, 2).Value = "RENAN"
$ws.Cells.Item(135, 3).Value = 46.76
$ws.Cells.Item(136, 1).Value = "001719494"
$ws.Cells.Item(136, 2).Value = "LUIS"
$ws.Cells.Item(136, 3).Value = 46.74
$ws.Cells.Item(137, 1).Value = "004277637"
$ws.Cells.Item(137, 2).Value = "LARA"
$ws.Cells.Item(137, 3).Value = 46.2
$ws.Cells.Item(138, 1).Value = "001731007"
$ws.Cells.Item(138, 2).Value = "GUILHERME"
$ws.Cells.Item(138, 3).Value = 44.59
$ws.Cells.Item(139, 1).Value = "004581652"
$ws.Cells.Item(139, 2).Value = "CINCO"
$ws.Cells.Item(139, 3).Value = 44.13
$ws.Cells.Item(140, 1).Value = "004805133"
$ws.Cells.Item(140, 2).Value = "PATRICIA"
$ws.Cells.Item(140, 3).Value = 41.48
$ws.Cells.Item(141, 1).Value = "004958578"
$ws.Cells.Item(141, 2).Value = "ASSAKO"
$ws.Cells.Item(141, 3).Value = 40.7
$ws.Cells.Item(142, 1).Value = "005165116"
$ws.Cells.Item(142, 2).Value = "ANA"
$ws.Cells.Item(142, 3).Value = 40.11
$ws.Cells.Item(143, 1).Value = "004998717"
$ws.Cells.Item(143, 2).Value = "GIOVANE"
$ws.Cells.Item(143, 3).Value = 40.08
$ws.Cells.Item(144, 1).Value = "004238164"
$ws.Cells.Item(144, 2).Value = "DANIELA"
$ws.Cells.Item(144, 3).Value = 38.3
$ws.Cells.Item(145, 1).Value = "004520100"
$ws.Cells.Item(145, 2).Value = "ALEXANDRE"
$ws.Cells.Item(145, 3).Value = 37.96
$ws.Cells.Item(146, 1).Value = "002401479"
$ws.Cells.Item(146, 2).Value = "JULIO"
$ws.Cells.Item(146, 3).Value = 37.84
$ws.Cells.Item(147, 1).Value = "004001621"
$ws.Cells.Item(147, 2).Value = "DANIELA"
$ws.Cells.Item(147, 3).Value = 37.58
$ws.Cells.Item(148, 1).Value = "005000656"
$ws.Cells.Item(148, 2).Value = "LUCIA"
$ws.Cells.Item(148, 3).Value = 35.88
$ws.Cells.Item(149, 1).Value = "004211922"
$ws.Cells.Item(149, 2).Value = "CARLOS"
$ws.Cells.Item(149, 3).Value = 34.71
$ws.Cells.Item(150, 1).Value = "004340984"
$ws.Cells.Item(150, 2).Value = "RENATA"
$ws.Cells.Item(150, 3).Value = 34
$ws.Cells.Item(151, 1).Value = "004994036"
$ws.Cells.Item(151, 2).Value = "BALTASAR"
$ws.Cells.Item(151, 3).Value = 33.73
$ws.Cells.Item(152, 1).Value = "004691225"
$ws.Cells.Item(152, 2).Value = "ANNA"
$ws.Cells.Item(152, 3).Value = 33.64
$ws.Cells.Item(153, 1).Value = "004472431"
$ws.Cells.Item(153, 2).Value = "LUIS"
$ws.Cells.Item(153, 3).Value = 33.08
$ws.Cells.Item(154, 1).Value = "004213943"
$ws.Cells.Item(154, 2).Value = "ELISA"
$ws.Cells.Item(154, 3).Value = 30.26
$ws.Cells.Item(155, 1).Value = "005018038"
$ws.Cells.Item(155, 2).Value = "ELAINE"
$ws.Cells.Item(155, 3).Value = 27.91
$ws.Cells.Item(156, 1).Value = "005105172"
$ws.Cells.Item(156, 2).Value = "VALDIVINO"
$ws.Cells.Item(156, 3).Value = 26.42
$ws.Cells.Item(157, 1).Value = "004377415"
$ws.Cells.Item(157, 2).Value = "ANGELA"
$ws.Cells.Item(157, 3).Value = 26.37
$ws.Cells.Item(158, 1).Value = "004240292"
$ws.Cells.Item(158, 2).Value = "MARCO"
$ws.Cells.Item(158, 3).Value = 24.3
$ws.Cells.Item(159, 1).Value = "004404724"
$ws.Cells.Item(159, 2).Value = "LEANDRO"
$ws.Cells.Item(159, 3).Value = 24.14
$ws.Cells.Item(160, 1).Value = "005173958"
$ws.Cells.Item(160, 2).Value = "VENIA"
$ws.Cells.Item(160, 3).Value = 23.27
$ws.Cells.Item(161, 1).Value = "004452507"
$ws.Cells.Item(161, 2).Value = "DANIELA"
$ws.Cells.Item(161, 3).Value = 22.74
$ws.Cells.Item(162, 1).Value = "005186167"
$ws.Cells.Item(162, 2).Value = "ANDREA"
$ws.Cells.Item(162, 3).Value = 21
$ws.Cells.Item(163, 1).Value = "004214604"
$ws.Cells.Item(163, 2).Value = "MARIA"
$ws.Cells.Item(163, 3).Value = 20.75
$ws.Cells.Item(164, 1).Value = "004920447"
$ws.Cells.Item(164, 2).Value = "MARILIA"
$ws.Cells.Item(164, 3).Value = 20.67
$ws.Cells.Item(165, 1).Value = "004458604"
$ws.Cells.Item(165, 2).Value = "FABIOLA"
$ws.Cells.Item(165, 3).Value = 20
$ws.Cells.Item(166, 1).Value = "004204255"
$ws.Cells.Item(166, 2).Value = "AMADO"
$ws.Cells.Item(166, 3).Value = 18.77
$ws.Cells.Item(167, 1).Value = "004368994"
$ws.Cells.Item(167, 2).Value = "CRISTINA"
$ws.Cells.Item(167, 3).Value = 18.56
$ws.Cells.Item(168, 1).Value = "004756968"
$ws.Cells.Item(168, 2).Value = "DANIELY"
$ws.Cells.Item(168, 3).Value = 18.08
$ws.Cells.Item(169, 1).Value = "004381194"
$ws.Cells.Item(169, 2).Value = "ALINNE"
$ws.Cells.Item(169, 3).Value = 17.71
$ws.Cells.Item(170, 1).Value = "001879977"
$ws.Cells.Item(170, 2).Value = "THAISSA"
$ws.Cells.Item(170, 3).Value = 17.14
$ws.Cells.Item(171, 1).Value = "005143579"
$ws.Cells.Item(171, 2).Value = "GABRIEL"
$ws.Cells.Item(171, 3).Value = 16.18
$ws.Cells.Item(172, 1).Value = "005169333"
$ws.Cells.Item(172, 2).Value = "EDUARDO"
$ws.Cells.Item(172, 3).Value = 16.12
$ws.Cells.Item(173, 1).Value = "004422594"
$ws.Cells.Item(173, 2).Value = "WANDIR"
$ws.Cells.Item(173, 3).Value = 14.67
$ws.Cells.Item(174, 1).Value = "004565108"
$ws.Cells.Item(174, 2).Value = "GUSTAVO"
$ws.Cells.Item(174, 3).Value = 14.51
$ws.Cells.Item(175, 1).Value = "000827730"
$ws.Cells.Item(175, 2).Value = "LUCIANA"
$ws.Cells.Item(175, 3).Value = 13.29
$ws.Cells.Item(176, 1).Value = "004752461"
$ws.Cells.Item(176, 2).Value = "SERGIO"
$ws.Cells.Item(176, 3).Value = 10.77
$ws.Cells.Item(177, 1).Value = "004216298"
$ws.Cells.Item(177, 2).Value = "FLORDELIZ"
$ws.Cells.Item(177, 3).Value = 9.75
$ws.Cells.Item(178, 1).Value = "004527606"
$ws.Cells.Item(178, 2).Value = "MARCIA"
$ws.Cells.Item(178, 3).Value = 9.52
$ws.Cells.Item(179, 1).Value = "004264780"
$ws.Cells.Item(179, 2).Value = "MARCELO"
$ws.Cells.Item(179, 3).Value = 8.99
$ws.Cells.Item(180, 1).Value = "004921978"
$ws.Cells.Item(180, 2).Value = "ELAINE"
$ws.Cells.Item(180, 3).Value = 8.08
$ws.Cells.Item(181, 1).Value = "004214460"
$ws.Cells.Item(181, 2).Value = "MARIA"
$ws.Cells.Item(181, 3).Value = 7.54
$ws.Cells.Item(182, 1).Value = "004693631"
$ws.Cells.Item(182, 2).Value = "NELY"
$ws.Cells.Item(182, 3).Value = 7.36
$ws.Cells.Item(183, 1).Value = "004530494"
$ws.Cells.Item(183, 2).Value = "ROSANGELA"
$ws.Cells.Item(183, 3).Value = 6.94
$ws.Cells.Item(184, 1).Value = "004854496"
$ws.Cells.Item(184, 2).Value = "JOSE"
$ws.Cells.Item(184, 3).Value = 6.64
$ws.Cells.Item(185, 1).Value = "004448501"
$ws.Cells.Item(185, 2).Value = "JOAO"
$ws.Cells.Item(185, 3).Value = 5.55
$ws.Cells.Item(186, 1).Value = "004216434"
$ws.Cells.Item(186, 2).Value = "JAIME"
$ws.Cells.Item(186, 3).Value = 4.83
$ws.Cells.Item(187, 1).Value = "005142624"
$ws.Cells.Item(187, 2).Value = "RODRIGO"
$ws.Cells.Item(187, 3).Value = 4.75
$ws.Cells.Item(188, 1).Value = "004239624"
$ws.Cells.Item(188, 2).Value = "NINA"
$ws.Cells.Item(188, 3).Value = 4.29
$ws.Cells.Item(189, 1).Value = "004848927"
$ws.Cells.Item(189, 2).Value = "ULDARICO"
$ws.Cells.Item(189, 3).Value = 3.62
$ws.Cells.Item(190, 1).Value = "005142661"
$ws.Cells.Item(190, 2).Value = "SABRINA"
$ws.Cells.Item(190, 3).Value = 3.6
$ws.Cells.Item(191, 1).Value = "004382374"
$ws.Cells.Item(191, 2).Value = "THEOMAR"
$ws.Cells.Item(191, 3).Value = 1.98
$ws.Cells.Item(192, 1).Value = "005341184"
$ws.Cells.Item(192, 2).Value = "BRENO"
$ws.Cells.Item(192, 3).Value = 1.85
$ws.Cells.Item(193, 1).Value = "004332783"
$ws.Cells.Item(193, 2).Value = "IRON"
$ws.Cells.Item(193, 3).Value = 1.73
$ws.Cells.Item(194, 1).Value = "004886366"
$ws.Cells.Item(194, 2).Value = "RENATO"
$ws.Cells.Item(194, 3).Value = 1.57
$ws.Cells.Item(195, 1).Value = "004308815"
$ws.Cells.Item(195, 2).Value = "ZELI"
$ws.Cells.Item(195, 3).Value = 1.25
$ws.Cells.Item(196, 1).Value = "005228239"
$ws.Cells.Item(196, 2).Value = "DEBORA"
$ws.Cells.Item(196, 3).Value = 0.85
$ws.Cells.Item(197, 1).Value = "004223502"
$ws.Cells.Item(197, 2).Value = "BRUNA"
$ws.Cells.Item(197, 3).Value = 0.78
$ws.Cells.Item(198, 1).Value = "004212581"
$ws.Cells.Item(198, 2).Value = "MARIA"
$ws.Cells.Item(198, 3).Value = 0.59
$ws.Cells.Item(199, 1).Value = "004453302"
$ws.Cells.Item(199, 2).Value = "ISABELLA"
$ws.Cells.Item(199, 3).Value = 0.39
$ws.Cells.Item(200, 1).Value = "004806286"
$ws.Cells.Item(200, 2).Value = "VERA"
$ws.Cells.Item(200, 3).Value = 0.19
$ws.Cells.Item(201, 1).Value = "004371857"
$ws.Cells.Item(201, 2).Value = "NAZARETH"
$ws.Cells.Item(201, 3).Value = 0.18
$ws.Cells.Item(202, 1).Value = "004332207"
$ws.Cells.Item(202, 2).Value = "IRACY"
$ws.Cells.Item(202, 3).Value = 0.16
$ws.Cells.Item(203, 1).Value = "004357159"
$ws.Cells.Item(203, 2).Value = "JOAO"
$ws.Cells.Item(203, 3).Value = 0.15
$ws.Cells.Item(204, 1).Value = "004320840"
$ws.Cells.Item(204, 2).Value = "NATALIA"
$ws.Cells.Item(204, 3).Value = 0.14
$ws.Cells.Item(205, 1).Value = "004466350"
$ws.Cells.Item(205, 2).Value = "RAQUEL"
$ws.Cells.Item(205, 3).Value = 0.11
$ws.Cells.Item(206, 1).Value = "005047946"
$ws.Cells.Item(206, 2).Value = "GABRIEL"
$ws.Cells.Item(206, 3).Value = 0.09
$ws.Cells.Item(207, 1).Value = "004589311"
$ws.Cells.Item(207, 2).Value = "CLARICE"
$ws.Cells.Item(207, 3).Value = 0.06
$ws.Cells.Item(208, 1).Value = "004321016"
$ws.Cells.Item(208, 2).Value = "JOAQUIM"
$ws.Cells.Item(208, 3).Value = 0.02
$ws.Cells.Item(209, 1).Value = "004850070"
$ws.Cells.Item(209, 2).Value = "RENATO"
$ws.Cells.Item(209, 3).Value = 0.02
$ws.Cells.Item(210, 1).Value = "002878817"
$ws.Cells.Item(210, 2).Value = "GUILHERME"
$ws.Cells.Item(210, 3).Value = 0.01
$ws.Cells.Item(211, 1).Value = "004400000"
$ws.Cells.Item(211, 2).Value = "VILMA"
$ws.Cells.Item(211, 3).Value = 0.01
$ws.Cells.Item(212, 1).Value = "004459875"
$ws.Cells.Item(212, 2).Value = "HELVECIO"
$ws.Cells.Item(212, 3).Value = 0.01
$ws.Cells.Item(213, 1).Value = "004612043"
$ws.Cells.Item(213, 2).Value = "YURI"
$ws.Cells.Item(213, 3).Value = 0.01
$ws.Cells.Item(214, 1).Value = "002823185"
$ws.Cells.Item(214, 2).Value = "SIMONE"
$ws.Cells.Item(214, 3).Value = -0.08
$ws.Cells.Item(215, 1).Value = "004335144"
$ws.Cells.Item(215, 2).Value = "EDMUNDO"
$ws.Cells.Item(215, 3).Value = -1.39
